# Apply the v3-spec metadata-field renames to the "assay" and "lanes"
# sheets, fix the one stray "10xV2" value, and restore the saved
# selection/view state to match the latest authoring session.

$wb = $excel.ActiveWorkbook

# --- "lanes" sheet: rename the read-file headers (done first so the
# new shared strings "r1"/"r2" land before the "assay" sheet's new
# strings, matching the authoring session's edit order) ---
$lanes = $wb.Worksheets.Item("lanes")
$lanes.Range("C1").Value = "r1"
$lanes.Range("D1").Value = "r2"
$lanes.Range("D2").Select()

# --- "assay" sheet: rename single_cell.* / seq.umi_barcode_* headers ---
$assay = $wb.Worksheets.Item("assay")
$assay.Range("B1").Value = "single_cell.cell_handling"
$assay.Range("C1").Value = "single_cell.barcode.offset"
$assay.Range("D1").Value = "single_cell.barcode.read"
$assay.Range("E1").Value = "single_cell.barcode.size"
$assay.Range("H1").Value = "seq.umi_barcode.offset"
$assay.Range("N1").Value = "seq.umi_barcode.size"
$assay.Range("O1").Value = "seq.umi_barcode.read"

# Data row: the rna.library_construction value used the old "10xV2"
# token; normalise it to match the "single_cell.cell_handling" value.
$assay.Range("L2").Value = "10x_v2"

# Widen the (now split) column L and restore the default view/selection.
# (ColumnWidth round-trips through a pixel conversion that adds 5/6 to
# the stored character width, so back the input off by that amount to
# land exactly on a stored width of 32.)
$assay.Columns.Item(12).ColumnWidth = 31.166666666666668
$assay.Application.ActiveWindow.ScrollColumn = 1
$assay.Range("C3").Select()
